$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0.5271737866561241
$ws.Range("G4").Value = 0.5802404724073328
$ws.Range("F5").Value = 0.2820299692918657
$ws.Range("G5").Value = 0.3163706874128024
$ws.Range("F6").Value = 0.1273308470964305
$ws.Range("G6").Value = 0.09794486453172768
